$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 200
$ws.Range("I33").Value = 200
$ws.Range("K33").Value = 200
$ws.Range("M33").Value = 29

$ws.Range("H70").Value = 101409.5
$ws.Range("J70").Value = 101409.5
$ws.Range("L70").Value = 304228.5
$ws.Range("N70").Value = -304768.5

$ws.Range("H73").Value = 101409.5
$ws.Range("J73").Value = 101409.5
$ws.Range("L73").Value = 304228.5
$ws.Range("N73").Value = -306100.5

$ws.Range("H88").Value = 1602.5834
$ws.Range("J88").Value = 1584.4546
$ws.Range("L88").Value = 1584.4546
$ws.Range("N88").Value = -2396.4546

$ws.Range("H91").Value = 1602.5834
$ws.Range("J91").Value = 1584.4546
$ws.Range("L91").Value = 1584.4546
$ws.Range("N91").Value = -4392.4546

$ws.Range("H112").Value = 3020.7
$ws.Range("I112").Value = 1799
$ws.Range("J112").Value = 3085
$ws.Range("K112").Value = 5397
$ws.Range("L112").Value = 9255
$ws.Range("M112").Value = -4289
$ws.Range("N112").Value = -11471

$ws.Range("H132").Value = 2891.6956
$ws.Range("I132").Value = 3096.6191
$ws.Range("K132").Value = 9289.8573
$ws.Range("M132").Value = -6759.8573

$ws.Range("H137").Value = 4320.7
$ws.Range("I137").Value = 4320.7
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 12962.1
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -10412.1
$ws.Range("N137").ClearContents()

$ws.Range("H141").Value = 2607.2354
$ws.Range("I141").Value = 2639.5625
$ws.Range("K141").Value = 7918.6875
$ws.Range("M141").Value = -2738.6875


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11868
$ws.Range("I32").Value = 11868
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 11868
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -11581
$ws.Range("N32").ClearContents()

$ws.Range("H63").Value = 5050.3335
$ws.Range("I63").Value = 2537.25
$ws.Range("J63").Value = 7060.8
$ws.Range("K63").Value = 2537.25
$ws.Range("L63").Value = 7060.8
$ws.Range("M63").Value = -1851.25
$ws.Range("N63").Value = -8432.799999999999

$ws.Range("H66").Value = 5050.3335
$ws.Range("I66").Value = 2537.25
$ws.Range("J66").Value = 7060.8
$ws.Range("K66").Value = 12686.25
$ws.Range("L66").Value = 35304
$ws.Range("M66").Value = -9254.25
$ws.Range("N66").Value = -42168

$ws.Range("H74").Value = 1699.8667
$ws.Range("I74").Value = 1699.8667
$ws.Range("K74").Value = 1699.8667
$ws.Range("M74").Value = -825.8667

$ws.Range("H77").Value = 1699.8667
$ws.Range("I77").Value = 1699.8667
$ws.Range("K77").Value = 8499.333500000001
$ws.Range("M77").Value = -4131.333500000001

$ws.Range("H97").Value = 1008.86664
$ws.Range("I97").Value = 1008.86664
$ws.Range("K97").Value = 1008.86664
$ws.Range("M97").Value = -512.86664

$ws.Range("H132").Value = 3213.6572
$ws.Range("I132").Value = 3205.2354
$ws.Range("K132").Value = 9615.706200000001
$ws.Range("M132").Value = -7085.706200000001


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2857.2856
$ws.Range("I99").Value = 2857.2856
$ws.Range("K99").Value = 2857.2856
$ws.Range("M99").Value = -1359.2856

$ws.Range("H134").Value = 77559.36
$ws.Range("I134").Value = 6602.385
$ws.Range("K134").Value = 19807.155
$ws.Range("M134").Value = -17272.155


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4714.75
$ws.Range("I58").Value = 4399.74
$ws.Range("J58").Value = 5430.6816
$ws.Range("K58").Value = 4399.74
$ws.Range("L58").Value = 5430.6816
$ws.Range("M58").Value = -4196.74
$ws.Range("N58").Value = -5836.6816

$ws.Range("H86").Value = 7131.5713
$ws.Range("I86").Value = 7039.3335
$ws.Range("J86").Value = 7297.6
$ws.Range("K86").Value = 7039.3335
$ws.Range("L86").Value = 7297.6
$ws.Range("M86").Value = -5916.3335
$ws.Range("N86").Value = -9543.6

$ws.Range("H89").Value = 7131.5713
$ws.Range("I89").Value = 7039.3335
$ws.Range("J89").Value = 7297.6
$ws.Range("K89").Value = 35196.6675
$ws.Range("L89").Value = 36488
$ws.Range("M89").Value = -29580.6675
$ws.Range("N89").Value = -47720

$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470

$ws.Range("H134").Value = 234735.25
$ws.Range("I134").Value = 2292.9
$ws.Range("K134").Value = 6878.700000000001
$ws.Range("M134").Value = -4343.700000000001

$ws.Range("H136").Value = 4714.75
$ws.Range("I136").Value = 4399.74
$ws.Range("J136").Value = 5430.6816
$ws.Range("K136").Value = 13199.22
$ws.Range("L136").Value = 16292.0448
$ws.Range("M136").Value = -10649.22
$ws.Range("N136").Value = -21392.0448


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 154.9375
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 154.9375
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 464.8125
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -810.8125

$ws.Range("H103").Value = 3821.5557
$ws.Range("J103").Value = 6736.8
$ws.Range("L103").Value = 20210.4
$ws.Range("N103").Value = -21968.4

$ws.Range("H124").Value = 944
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H131").Value = 3206
$ws.Range("I131").Value = 1311.4445
$ws.Range("K131").Value = 3934.3335
$ws.Range("M131").Value = 1105.6665


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 834.1579
$ws.Range("I102").Value = 605.6875
$ws.Range("K102").Value = 605.6875
$ws.Range("M102").Value = 1016.3125

$ws.Range("H126").Value = 33336462
$ws.Range("I126").Value = 55558160
$ws.Range("J126").Value = 3915.8333
$ws.Range("K126").Value = 166674480
$ws.Range("L126").Value = 11747.4999
$ws.Range("M126").Value = -166672010
$ws.Range("N126").Value = -16687.4999

$ws.Range("H132").Value = 119722.22
$ws.Range("I132").Value = 13700.4
$ws.Range("K132").Value = 41101.2
$ws.Range("M132").Value = -38571.2

$ws.Range("H135").Value = 500149980
$ws.Range("J135").Value = 500149980
$ws.Range("L135").Value = 500149980
$ws.Range("N135").Value = -500160120

$ws.Range("H137").Value = 49999
$ws.Range("I137").Value = 49999
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 49999
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -44899
$ws.Range("N137").ClearContents()

$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

$ws.Range("H140").Value = 200000
$ws.Range("J140").Value = 200000
$ws.Range("L140").Value = 200000
$ws.Range("N140").Value = -210360


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 200001220
$ws.Range("I16").Value = 333334200
$ws.Range("K16").Value = 333334200
$ws.Range("M16").Value = -333334030

$ws.Range("H46").Value = 3865.3845
$ws.Range("I46").Value = 3500
$ws.Range("J46").Value = 4687.5
$ws.Range("K46").Value = 3500
$ws.Range("L46").Value = 4687.5
$ws.Range("M46").Value = -3312
$ws.Range("N46").Value = -5063.5

$ws.Range("H55").Value = 1328.3334
$ws.Range("I55").Value = 301.8889
$ws.Range("K55").Value = 301.8889
$ws.Range("M55").Value = -128.8889

$ws.Range("H61").Value = 3309.9285
$ws.Range("I61").Value = 2867.16
$ws.Range("K61").Value = 2867.16
$ws.Range("M61").Value = -2665.16

$ws.Range("H93").Value = 90911480
$ws.Range("I93").Value = 111113210
$ws.Range("K93").Value = 111113210
$ws.Range("M93").Value = -111111962

$ws.Range("H113").Value = 3309.9285
$ws.Range("I113").Value = 2867.16
$ws.Range("K113").Value = 2867.16
$ws.Range("M113").Value = -697.1599999999999

$ws.Range("H122").Value = 1822071.1
$ws.Range("I122").Value = 1431285.9
$ws.Range("J122").Value = 2505945.2
$ws.Range("K122").Value = 4293857.699999999
$ws.Range("L122").Value = 7517835.600000001
$ws.Range("M122").Value = -4291407.699999999
$ws.Range("N122").Value = -7522735.600000001

$ws.Range("H132").Value = 7327.1875
$ws.Range("I132").Value = 6168.125
$ws.Range("J132").Value = 8486.25
$ws.Range("K132").Value = 18504.375
$ws.Range("L132").Value = 25458.75
$ws.Range("M132").Value = -15974.375
$ws.Range("N132").Value = -30518.75


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 77085.5
$ws.Range("I62").Value = 147028.14
$ws.Range("K62").Value = 147028.14
$ws.Range("M62").Value = -146404.14

$ws.Range("H65").Value = 77085.5
$ws.Range("I65").Value = 147028.14
$ws.Range("K65").Value = 735140.7000000001
$ws.Range("M65").Value = -732020.7000000001

$ws.Range("H81").Value = 1742.4706
$ws.Range("J81").Value = 2374
$ws.Range("L81").Value = 4748
$ws.Range("N81").Value = -6870

$ws.Range("H84").Value = 1742.4706
$ws.Range("J84").Value = 2374
$ws.Range("L84").Value = 23740
$ws.Range("N84").Value = -34348

$ws.Range("H100").Value = 755.1
$ws.Range("I100").Value = 864
$ws.Range("J100").Value = 501
$ws.Range("K100").Value = 1728
$ws.Range("L100").Value = 1002
$ws.Range("M100").Value = -1187
$ws.Range("N100").Value = -2084

$ws.Range("H122").Value = 50002170
$ws.Range("I122").Value = 62501412
$ws.Range("J122").Value = 5200.75
$ws.Range("K122").Value = 187504236
$ws.Range("L122").Value = 15602.25
$ws.Range("M122").Value = -187501786
$ws.Range("N122").Value = -20502.25

$ws.Range("H132").Value = 36918.035
$ws.Range("I132").Value = 2437.0435
$ws.Range("J132").Value = 150212.72
$ws.Range("K132").Value = 7311.130500000001
$ws.Range("L132").Value = 450638.16
$ws.Range("M132").Value = -4781.130500000001
$ws.Range("N132").Value = -455698.16

